$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

# Update the date in A1 (was 2024-04-24 / serial 45406, now 2024-05-24 / serial 45436)
$ws.Range("A1").Value = Get-Date -Year 2024 -Month 5 -Day 24 -Hour 0 -Minute 0 -Second 0

# Update the prices in column D for rows 33-35
$ws.Range("D33").Value = 1507.712
$ws.Range("D34").Value = 2542.857
$ws.Range("D35").Value = 2799.394
